$wb = $excel.ActiveWorkbook

# --- Sheet "Đơn phụ phẫu 1": insert a new detail row before the "Tổng" row ---
$ws = $wb.Worksheets.Item("Đơn phụ phẫu 1")

# Push the existing "Tổng" row (currently row 15) down to row 16, making room
# for the new HD-LUXURY entry at row 15.
$ws.Rows.Item(15).Insert()

$ws.Cells.Item(15, 1).Value = "HD-LUXURY"
$ws.Cells.Item(15, 2).Value = 692
# The service-date column holds plain text like "08-27-2024" (not a real
# date), so prefix with an apostrophe to keep Excel from reinterpreting it
# as a date serial number.
$ws.Cells.Item(15, 3).Value = "'08-27-2024"
$ws.Cells.Item(15, 4).Value = "CẦN THƠ"
$ws.Cells.Item(15, 5).Value = "Phan Thị Diễm My"
$ws.Cells.Item(15, 6).Value = "Cá nhân"
$ws.Cells.Item(15, 7).Value = "Nâng mũi"
$ws.Cells.Item(15, 8).Value = "Lâm Hoàng Phú"
$ws.Cells.Item(15, 9).Value = 100000

# Update the "Tổng" (Total) row, now on row 16: one more order, +100000 total.
$ws.Cells.Item(16, 2).Value = 14
$ws.Cells.Item(16, 9).Value = 1000000

# --- Sheet "Lương": roll the new Phụ phẫu 1 commission into the salary summary ---
$wsLuong = $wb.Worksheets.Item("Lương")

$wsLuong.Cells.Item(8, 2).Value = 1000000
$wsLuong.Cells.Item(34, 2).Value = 923571.4285714286
$wsLuong.Cells.Item(37, 2).Value = 1023571.428571429
